# installer will add shortcut on desktop
#
# 1. Update existing Title text on the Active sheet (Id=22) to mention
#    that the msi-installed program needs to run again after installation.
# 2. Insert a new task row on the Active sheet for "better icon"
#    (Id=38, Status=Todo, Category=Feature, Created=4/11/2018).
# 3. Bump the "Max Id" tracker on the Config sheet from 37 to 38 to
#    match the newly used Id.

$wb = $excel.ActiveWorkbook

$wsActive = $wb.Worksheets.Item("Active")

# 1. Update row 2 (Id 22) title text.
$wsActive.Cells.Item(2, 2).Value = "get the msi installed program to run again - it won't run after installation"

# 2. Insert a new row at position 6 (pushes existing rows 6-9 down to 7-10)
#    and populate it with the new "better icon" task.
$wsActive.Rows.Item(6).Insert()

$wsActive.Cells.Item(6, 1).Value = 38
$wsActive.Cells.Item(6, 2).Value = "better icon"
$wsActive.Cells.Item(6, 3).Value = "Todo"
$wsActive.Cells.Item(6, 4).Value = "Feature"
# Force text (not an auto-converted date serial) to match the existing
# "Created" column entries, which are all plain text like "4/11/2018".
$wsActive.Cells.Item(6, 5).NumberFormat = "@"
$wsActive.Cells.Item(6, 5).Value = "4/11/2018"

# Match the default (unbolded) style used by the rest of the data rows.
$wsActive.Range("A6:E6").Style = $wsActive.Range("A5:E5").Style

# 3. Bump the Max Id counter on the Config sheet.
$wsConfig = $wb.Worksheets.Item("Config")
$wsConfig.Cells.Item(2, 6).Value = 38
